$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.256.42"
$ws.Range("E2").Value = "  -5.28%  "
$ws.Range("D3").Value = "2.898.50"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "2.893.71"
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "3.361.99"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "2.883.73"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("D19").Value = "57.164.15"
$ws.Range("E19").Value = "  -5.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "405.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0980"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.913"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0641"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.106"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0337"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "363.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").Value = "2.605.10"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.229"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "
